$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = 30
$ws.Range("C14").Value = 0.001
$ws.Range("D14").Value = 0.003
$ws.Range("E14").Value = "Regular"
$ws.Range("F14").Value = "<function relu at 0x11ad159d8>"
$ws.Range("G14").Value = 0.9657999873161316
$ws.Range("H14").Value = 0.2011000066995621
$ws.Range("I14").Value = 0.1811999976634979
$ws.Range("J14").Value = 0.1451183259487152
$ws.Range("K14").Value = 5.327450752258301
$ws.Range("L14").Value = 0.2011000066995621
$ws.Range("M14").Value = "logs/results_71.log"
$ws.Range("N14").Value = "weights/model_71.ckpt"
$ws.Range("O14").Value = "tb/71/non_robust"
$ws.Range("P14").Value = "(5.475276, 12.719564, 18.869154, 27.198263, 26.215324, 22.13533, 16.86695)"
$ws.Range("Q14").Value = "(78.487724, 11.971958, 8.304008, 4.5321946, 2.381914, 1.7097418, 1.7520251, 2.2926486)"

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = 30
$ws.Range("C15").Value = 0.005
$ws.Range("D15").Value = 0.003
$ws.Range("E15").Value = "Regular"
$ws.Range("F15").Value = "<function relu at 0x12008f9d8>"
$ws.Range("G15").Value = 0.9585999846458435
$ws.Range("H15").Value = 0.2369000017642975
$ws.Range("I15").Value = 0.03700000047683716
$ws.Range("J15").Value = 0.1429557055234909
$ws.Range("K15").Value = 4.122503280639648
$ws.Range("L15").Value = 0.2369000017642975
$ws.Range("M15").Value = "logs/results_72.log"
$ws.Range("N15").Value = "weights/model_72.ckpt"
$ws.Range("O15").Value = "tb/72/non_robust"
$ws.Range("P15").Value = "(2.139943, 3.496155, 3.843952, 5.1138744, 6.2203045, 6.9547405, 8.452685)"
$ws.Range("Q15").Value = "(29.190838, 5.508862, 3.9576876, 1.9097207, 1.488211, 1.3313731, 1.552425, 1.4604229)"

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = 70
$ws.Range("C16").Value = 0.008
$ws.Range("D16").Value = 0.003
$ws.Range("E16").Value = "Regular"
$ws.Range("F16").Value = "<function relu at 0x1240199d8>"
$ws.Range("G16").Value = 0.9384999871253967
$ws.Range("H16").Value = 0.2345000058412552
$ws.Range("I16").Value = 0.07829999923706055
$ws.Range("J16").Value = 0.2154168486595154
$ws.Range("K16").Value = 4.516582012176514
$ws.Range("L16").Value = 0.2345000058412552
$ws.Range("M16").Value = "logs/results_74.log"
$ws.Range("N16").Value = "weights/model_74.ckpt"
$ws.Range("O16").Value = "tb/74/non_robust"
$ws.Range("P16").Value = "(1.7076006, 2.5921235, 2.5732467, 3.4120953, 4.31576, 5.6463065, 6.7571893)"
$ws.Range("Q16").Value = "(23.436884, 4.6997, 2.8007762, 1.6223694, 1.4685857, 1.5694603, 1.4128212, 1.4283097)"

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = 30
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 0.003
$ws.Range("E17").Value = "FGSM"
$ws.Range("F17").Value = "<function relu at 0x11951f9d8>"
$ws.Range("G17").Value = 0.9517999887466431
$ws.Range("H17").Value = 0.8098999857902527
$ws.Range("I17").Value = 0.4519999921321869
$ws.Range("J17").Value = 0.1830078810453415
$ws.Range("K17").Value = 0.6702156066894531
$ws.Range("L17").Value = 0.8098999857902527
$ws.Range("M17").Value = "logs/results_75.log"
$ws.Range("N17").Value = "weights/model_75.ckpt"
$ws.Range("O17").Value = "tb/75/robust"
$ws.Range("P17").Value = "(6.767173, 17.473831, 38.257633, 49.1879, 27.285254, 10.4143715, 4.247443)"
$ws.Range("Q17").Value = "(131.64798, 17.081043, 17.57381, 15.574147, 13.246862, 13.433027, 11.568809, 11.362484)"
